# HP: ExcelToStringArray extended with write methods
# Append 5 new rows (21-25) of shared-string values to Sheet1, mirroring
# the "ExcelToStringArray" write-back behaviour: a new F-column item per
# row, mirrored into G starting the following row, trailing off at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F21").Value = "nitem100"

$ws.Range("F22").Value = "nitem101"
$ws.Range("G22").Value = "nitem101"

$ws.Range("F23").Value = "nitem102"
$ws.Range("G23").Value = "nitem102"

$ws.Range("F24").Value = "nitem103"
$ws.Range("G24").Value = "nitem103"

$ws.Range("G25").Value = "nitem104"
